$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove now-unused columns H:J (duplicate "combined" metrics no longer tracked)
$ws.Range("H1:J25").Delete()

# Row A/C order changes + refreshed values for the remaining quarters
$ws.Range("A2").Value = "2016年A"
$ws.Range("B2").Value = 78.11610820734271
$ws.Range("C2").Value = 78.11610820734271
$ws.Range("D2").Value = -16.4176283780316
$ws.Range("E2").Value = -16.4176283780316
$ws.Range("F2").Value = 38.3015201706892
$ws.Range("G2").Value = 38.3015201706892

$ws.Range("A3").Value = "2016年C"
$ws.Range("B3").Value = 58.3545560060529
$ws.Range("C3").Value = 63.6053102446392
$ws.Range("D3").Value = -10.7452617742999
$ws.Range("E3").Value = -10.558411982935
$ws.Range("F3").Value = 52.3907057682471
$ws.Range("G3").Value = 46.9531017382959

$ws.Range("A4").Value = "2016年B"
$ws.Range("B4").Value = 55.9217839058424
$ws.Range("C4").Value = 66.49133029705099
$ws.Range("D4").Value = -5.03498468015939
$ws.Range("E4").Value = -10.4557120240023
$ws.Range("F4").Value = 49.113200774317
$ws.Range("G4").Value = 43.9643817269513

$ws.Range("A5").Value = "2016年D"
$ws.Range("B5").Value = 72.01020425594039
$ws.Range("C5").Value = 65.9691414801613
$ws.Range("D5").Value = -14.5144773657391
$ws.Range("E5").Value = -11.6710340944175
$ws.Range("F5").Value = 42.5042731097988
$ws.Range("G5").Value = 45.7018926142561

$ws.Range("A6").Value = "2017年A"
$ws.Range("B6").Value = 76.1797512778989
$ws.Range("C6").Value = 76.1797512778989
$ws.Range("D6").Value = -4.54979825579324
$ws.Range("E6").Value = -4.54979825579324
$ws.Range("F6").Value = 28.3700469778939
$ws.Range("G6").Value = 28.3700469778939

$ws.Range("A7").Value = "2017年C"
$ws.Range("B7").Value = 60.1319431891811
$ws.Range("C7").Value = 61.5369662964737
$ws.Range("D7").Value = 1.47138532318578
$ws.Range("E7").Value = -0.723654200654494
$ws.Range("F7").Value = 38.396671487633
$ws.Range("G7").Value = 39.1866879041806

$ws.Range("A8").Value = "2017年B"
$ws.Range("B8").Value = 49.725724151169
$ws.Range("C8").Value = 62.3010344657317
$ws.Range("D8").Value = 0.46788842204143
$ws.Range("E8").Value = -1.91734263305418
$ws.Range("F8").Value = 49.8063874267897
$ws.Range("G8").Value = 39.6163081673223

$ws.Range("A9").Value = "2017年D"
$ws.Range("B9").Value = 42.4141360008699
$ws.Range("C9").Value = 55.8676814771835
$ws.Range("D9").Value = 17.4408282280971
$ws.Range("E9").Value = 4.6615122674735
$ws.Range("F9").Value = 40.1450357710331
$ws.Range("G9").Value = 39.4708062553431

$ws.Range("A10").Value = "2018年A"
$ws.Range("B10").Value = 68.5217116883003
$ws.Range("C10").Value = 68.5217116883003
$ws.Range("D10").Value = -16.3229161572722
$ws.Range("E10").Value = -16.3229161572722
$ws.Range("F10").Value = 47.8012044689721
$ws.Range("G10").Value = 47.8012044689721

$ws.Range("A11").Value = "2018年C"
$ws.Range("B11").Value = 62.9067460152289
$ws.Range("C11").Value = 66.5262791168975
$ws.Range("D11").Value = -7.65084708135837
$ws.Range("E11").Value = -10.71080233729
$ws.Range("F11").Value = 44.7441010661294
$ws.Range("G11").Value = 44.1845232203926

$ws.Range("A12").Value = "2018年B"
$ws.Range("B12").Value = 68.2836608024575
$ws.Range("C12").Value = 68.3983276493243
$ws.Range("D12").Value = -8.548636771870839
$ws.Range("E12").Value = -12.293433186423
$ws.Range("F12").Value = 40.2649759694132
$ws.Range("G12").Value = 43.8951055370988

$ws.Range("A13").Value = "2018年D"
$ws.Range("B13").Value = 57.4559617321919
$ws.Range("C13").Value = 63.9848715382958
$ws.Range("D13").Value = 1.87017980373149
$ws.Range("E13").Value = -7.18574340017295
$ws.Range("F13").Value = 40.6738584640767
$ws.Range("G13").Value = 43.2008718618773

$ws.Range("A14").Value = "2019年A"
$ws.Range("B14").Value = 66.0908618710563
$ws.Range("C14").Value = 66.0908618710563
$ws.Range("D14").Value = 20.5112145564399
$ws.Range("E14").Value = 20.5112145564399
$ws.Range("F14").Value = 13.3979235725038
$ws.Range("G14").Value = 13.3979235725038

$ws.Range("A15").Value = "2019年C"
$ws.Range("B15").Value = 59.1965948990418
$ws.Range("C15").Value = 59.0671637862615
$ws.Range("D15").Value = 17.4756854159693
$ws.Range("E15").Value = 17.2878594230108
$ws.Range("F15").Value = 23.3277196849889
$ws.Range("G15").Value = 23.6449767907276

$ws.Range("A16").Value = "2019年B"
$ws.Range("B16").Value = 52.2899088581393
$ws.Range("C16").Value = 58.9978891316221
$ws.Range("D16").Value = 14.0438601086188
$ws.Range("E16").Value = 17.1873304194735
$ws.Range("F16").Value = 33.6662310332419
$ws.Range("G16").Value = 23.8147804489044

$ws.Range("A17").Value = "2019年D"
$ws.Range("B17").Value = 57.2041753518251
$ws.Range("C17").Value = 58.5625527315661
$ws.Range("D17").Value = -0.19300343459329
$ws.Range("E17").Value = 12.552974074016
$ws.Range("F17").Value = 42.988828082768
$ws.Range("G17").Value = 28.8844731944176

$ws.Range("A18").Value = "2020年A"
$ws.Range("B18").Value = 58.4146963512201
$ws.Range("C18").Value = 58.4146963512201
$ws.Range("D18").Value = 17.5443429254068
$ws.Range("E18").Value = 17.5443429254068
$ws.Range("F18").Value = 24.0409607233731
$ws.Range("G18").Value = 24.0409607233731

$ws.Range("A19").Value = "2020年C"
$ws.Range("B19").Value = 35.5188093814298
$ws.Range("C19").Value = -302.404863073248
$ws.Range("D19").Value = 23.9591918769264
$ws.Range("E19").Value = 34.4474976043345
$ws.Range("F19").Value = 40.5219987416439
$ws.Range("G19").Value = 367.957365468917

$ws.Range("A20").Value = "2020年B"
$ws.Range("B20").Value = -68.4341273757271
$ws.Range("C20").Value = 155.550639361443
$ws.Range("D20").Value = 14.0323340435744
$ws.Range("E20").Value = 20.2337040449976
$ws.Range("F20").Value = 154.401793332153
$ws.Range("G20").Value = -75.7843434064407

$ws.Range("A21").Value = "2020年D"
$ws.Range("B21").Value = 44.8113002027044
$ws.Range("C21").Value = -6.84224030926098
$ws.Range("D21").Value = 23.7361927159421
$ws.Range("E21").Value = 25.3296577659804
$ws.Range("F21").Value = 31.4525070813535
$ws.Range("G21").Value = 81.5125825432814

$ws.Range("A22").Value = "2021年A"
$ws.Range("B22").Value = 47.1910071757524
$ws.Range("C22").Value = 47.1910071757524
$ws.Range("D22").Value = 25.7616412004057
$ws.Range("E22").Value = 25.7616412004057
$ws.Range("F22").Value = 27.0473516238418
$ws.Range("G22").Value = 27.0473516238418

$ws.Range("A23").Value = "2021年C"
$ws.Range("B23").Value = 65.251877166353
$ws.Range("C23").Value = 56.0709907359548
$ws.Range("D23").Value = 20.6457427381183
$ws.Range("E23").Value = 20.9147698132891
$ws.Range("F23").Value = 14.1023800955286
$ws.Range("G23").Value = 23.0142394507561

$ws.Range("A24").Value = "2021年B"
$ws.Range("B24").Value = 66.52705768914331
$ws.Range("C24").Value = 53.9930560068381
$ws.Range("D24").Value = 12.1566249077651
$ws.Range("E24").Value = 20.9756594350928
$ws.Range("F24").Value = 21.3163174030915
$ws.Range("G24").Value = 25.031284558069

$ws.Range("A25").Value = "2021年D"
$ws.Range("B25").Value = 70.6126698254599
$ws.Range("C25").Value = 58.261396957157
$ws.Range("D25").Value = 27.6196866029018
$ws.Range("E25").Value = 21.9247282478147
$ws.Range("F25").Value = 1.76764357163774
$ws.Range("G25").Value = 19.8138747950282

